$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.671.83'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.603.86'
$ws.Range('E3').Value = '  -0.03%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.34'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.84'
$ws.Range('E8').Value = '  +8.15%  '
$ws.Range('E9').Value = '  +2.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0606'
$ws.Range('E10').Value = '  +1.59%  '
$ws.Range('E11').Value = '  -0.50%  '
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.611.80'
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.556'
$ws.Range('E14').Value = '  +3.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.675.04'
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.80'
$ws.Range('E16').Value = '  +1.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.24'
$ws.Range('E17').Value = '  +1.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '240.78'
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.09'
$ws.Range('E19').Value = '  +6.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0702'
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('E23').Value = '  +3.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.10'
$ws.Range('E24').Value = '  +1.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.03'
$ws.Range('E25').Value = '  +1.02%  '
$ws.Range('E26').Value = '  +1.32%  '
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.48'
$ws.Range('E28').Value = '  +1.90%  '
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0479'
$ws.Range('E30').Value = '  +2.05%  '
$ws.Range('E31').Value = '  +0.53%  '
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('E33').Value = '  +2.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.424.47'
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('E35').Value = '  +3.92%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.84'
$ws.Range('E37').Value = '  +1.53%  '
$ws.Range('E38').Value = '  +0.44%  '
$ws.Range('E39').Value = '  +2.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.550'
$ws.Range('E40').Value = '  +3.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '55.33'
$ws.Range('E41').Value = '  +3.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0495'
$ws.Range('E42').Value = '  +4.87%  '
$ws.Range('E43').Value = '  +1.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.819'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '67.58'
$ws.Range('E46').Value = '  +2.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.996'
$ws.Range('E47').Value = '  +18.94%  '
$ws.Range('E48').Value = '  +2.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.741.29'
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('E50').Value = '  -1.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '86.68'
$ws.Range('E51').Value = '  -0.02%  '
